# Updates cryptos list figures (price/volume) to the latest scrape.
# Three rows (43-45) were also reordered: FirstDigitalUSD, ApeXProtocol, PEPE.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.616.24"
$ws.Range("E2").Value = "  +0.20%  "
# Row 3
$ws.Range("D3").Value = "3.774.68"
$ws.Range("E3").Value = "  -1.43%  "
# Row 4
$ws.Range("E4").Value = "  +0.07%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "434.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.61%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.39%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.620"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.21%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.07%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.732"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.21%  "
# Row 10
$ws.Range("E10").Value = "  -9.61%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000313"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -14.76%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.87%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.41%  "
# Row 14
$ws.Range("D14").Value = "4.393.16"
$ws.Range("E14").Value = "  -0.77%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.32%  "
# Row 17
$ws.Range("D17").Value = "3.789.14"
$ws.Range("E17").Value = "  -0.65%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.33%  "
# Row 19
$ws.Range("E19").Value = "  +6.79%  "
# Row 20
$ws.Range("D20").Value = "66.730.03"
$ws.Range("E20").Value = "  -0.13%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "407.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.95%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.51%  "
# Row 23
$ws.Range("E23").Value = "  +6.70%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.10%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "36.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.46%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.43%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +38.15%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.95%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.48%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "725.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.08%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.91%  "
# Row 32
$ws.Range("E32").Value = "  +9.52%  "
# Row 33
$ws.Range("E33").Value = "  +2.59%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.29%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.22%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.154"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.31%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +25.19%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "56.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.17%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0475"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.00%  "
# Row 40
$ws.Range("E40").Value = "  +41.68%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.08%  "
# Row 42
$ws.Range("E42").Value = "  +3.30%  "
# Row 43
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.20%  "
# Row 44
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.35%  "
# Row 45
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").Value = "0.0₃0667"
$ws.Range("E45").Value = "  -16.05%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.329"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +11.96%  "
# Row 47
$ws.Range("E47").Value = "  +5.56%  "
# Row 48
$ws.Range("E48").Value = "  +0.31%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.06%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.37%  "
# Row 51
$ws.Range("E51").Value = "  +0.89%  "
